$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet has two additional rows of data imported from an external
# source; they duplicate the existing "Override_Product_Information"
# record (row 3).
$ws.Range("A3:F3").Copy()
$ws.Range("A4:F4").PasteSpecial()
$ws.Range("A3:F3").Copy()
$ws.Range("A5:F5").PasteSpecial()
$excel.CutCopyMode = $false

# Reflect final cursor/selection position left behind after the import.
$ws.Range("E11").Select()
